$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.685507
$ws.Range("H2").Value = 11.056521
$ws.Range("I2").Value = 0.3585631737883472
$ws.Range("J2").Value = 0.3585631737883472
$ws.Range("Q2").Value = 0.1264251751233333
$ws.Range("R2").Value = 1.13782657611
$ws.Range("S2").Value = 0.3585631737883472
$ws.Range("T2").Value = 0.3585631737883472

# Row 3
$ws.Range("I3").Value = 0.009647184430711629
$ws.Range("J3").Value = 0.009647184430711629
$ws.Range("S3").Value = 0.009647184430711629
$ws.Range("T3").Value = 0.009647184430711629

# Row 4
$ws.Range("G4").Value = 6.493877
$ws.Range("H4").Value = 19.481631
$ws.Range("I4").Value = 0.6317896417809412
$ws.Range("J4").Value = 0.6317896417809411
$ws.Range("Q4").Value = 0.2227616273566667
$ws.Range("R4").Value = 2.00485464621
$ws.Range("S4").Value = 0.6317896417809412
$ws.Range("T4").Value = 0.6317896417809411
